$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Range("C$r")
    if ($cell.Value2 -eq 3462) {
        $cell.Value2 = 3463
    }
}
